$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# commit tabela de controle - add new row of data (nome / od / posto)
$ws.Range("A4").Value = "DEIVID ROMULO DA SILVA VICENTE"
$ws.Range("B4").Value = "Ordenador de Despesa Substituto"
$ws.Range("C4").Value = "Segundo-Sargento (FR)"

# B4 picks up the same cell formatting already used by B3 (same text, column)
$ws.Range("B3").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C7").Select() | Out-Null
